$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 6
    6  = 2
    7  = 5
    8  = 6
    9  = 6
    10 = 8
    11 = 3
    12 = 6
    13 = 0
    14 = 5
    15 = 6
    16 = 13
    17 = 5
    18 = 7
    19 = 3
    20 = 4
    21 = 5
    22 = 5
    23 = 6
    24 = 6
    25 = 6
    26 = 7
    27 = 6
    28 = 5
    29 = 3
    30 = 4
    31 = 7
    32 = 3
    33 = 9
    34 = 7
    35 = 6
    36 = 2
    37 = 5
    38 = 9
    39 = 3
    40 = 4
    41 = 1
    42 = 4
    43 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
